$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("BDPbES")

# "About" sheet: record the state this data applies to and refresh the date stamp.
$wsAbout.Range("B1").Value = "California"
$wsAbout.Range("C1").Value = "9/24/2021"

# "BDPbES" sheet: onshore wind and solar PV now get top dispatch priority (1 instead of 2).
$wsData.Range("B6").Value = 1
$wsData.Range("B7").Value = 1

# Make the BDPbES sheet the active/selected tab, with B8 selected.
$wsData.Activate()
$wsData.Range("B8").Select()
